# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Row => new value for column F
$updates = @{
    2  = 623
    4  = 79
    5  = 12955
    8  = 511
    12 = 13730
    13 = 14230
    15 = 171
    20 = 2
    25 = 936
    26 = 5312
    28 = 289
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
